$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "EC" (Estado de Cuenta) database: the "Periodo Mora" column
# (E16:E22) is re-sorted into ascending chronological order (1911 .. 2005)
# instead of the previous descending order (2005 .. 1911).
$ws.Range("E16").Value = "1911"
$ws.Range("E17").Value = "1912"
$ws.Range("E18").Value = "2001"
$ws.Range("E19").Value = "2002"
$ws.Range("E20").Value = "2003"
$ws.Range("E21").Value = "2004"
$ws.Range("E22").Value = "2005"

# The "Valor Mora" figures follow their period row, so the two rows whose
# period label moved to the opposite end of the list (now 1911 / 2005)
# swap their values accordingly.
$ws.Range("F16").Value = 33125
$ws.Range("F22").Value = 20979
